# Quarterly cash-flow sheet: roll the reporting window forward by one
# quarter (drop the oldest quarter column, add the new Q3 1401/09 column)
# and bump the copyright year. "add market cap to price" commit -- the
# visible effect in this workbook is adding the newest quarter's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copyright year bump -------------------------------------------------
$ws.Range("B3").Value2 = "Copyright @2015 - 2023"

# --- Header row 8: quarter labels shift left by one, newest quarter added
$ws.Range("D8").Value2 = "فصل سوم منتهی به 1400/09"
$ws.Range("E8").Value2 = "فصل چهارم منتهی به 1400/12"
$ws.Range("F8").Value2 = "فصل اول منتهی به 1401/03"
$ws.Range("G8").Value2 = "فصل دوم منتهی به 1401/06"
$ws.Range("H8").Value2 = "فصل سوم منتهی به 1401/09"

# --- Header row 9: publish dates shift left by one, newest date added ----
# These look like ISO dates, so Excel's type-inference would silently turn
# the "clean" ones (no "(n)" suffix) into date serials via .Value/.Value2.
# Force text by going through .Formula with a leading apostrophe instead,
# same as typing '1401-04-28 into the cell.
$ws.Range("D9").Formula = "'1401-10-29 (2)"
$ws.Range("E9").Formula = "'1401-10-29 (6)"
$ws.Range("F9").Formula = "'1401-04-28"
$ws.Range("G9").Formula = "'1401-08-25 (2)"
$ws.Range("H9").Formula = "'1401-10-29"

# --- Data rows: shift existing D:H values left by one column and fill ----
# the freed-up column H with the new quarter's figures.
# NOTE: `.Value` reads are unreliable in this host for chained property
# access, so use `.Value2` (also fine for writes) throughout.
function Shift-Row($row, $newH) {
    $e = $ws.Range("E$row").Value2
    $f = $ws.Range("F$row").Value2
    $g = $ws.Range("G$row").Value2
    $h = $ws.Range("H$row").Value2
    $ws.Range("D$row").Value2 = $e
    $ws.Range("E$row").Value2 = $f
    $ws.Range("F$row").Value2 = $g
    $ws.Range("G$row").Value2 = $h
    $ws.Range("H$row").Value2 = $newH
}

Shift-Row 12 -852487
Shift-Row 13 -51799
Shift-Row 14 -904286
Shift-Row 16 0
Shift-Row 17 -40942
Shift-Row 18 0
Shift-Row 19 0
Shift-Row 20 0
Shift-Row 21 0
Shift-Row 22 0
Shift-Row 23 0
Shift-Row 24 0
Shift-Row 25 0
Shift-Row 26 0
Shift-Row 27 0
Shift-Row 28 0
Shift-Row 29 0
Shift-Row 30 0
Shift-Row 31 322
Shift-Row 32 -40620
Shift-Row 33 -944906
Shift-Row 35 0
Shift-Row 37 0
Shift-Row 38 0
Shift-Row 39 900000
Shift-Row 40 -243520
Shift-Row 41 -27835
Shift-Row 42 0
Shift-Row 43 0
Shift-Row 44 0
Shift-Row 45 0
Shift-Row 46 0
Shift-Row 47 0
Shift-Row 48 0
Shift-Row 49 0
Shift-Row 50 -191374
Shift-Row 51 437271
Shift-Row 52 -507635
Shift-Row 53 548992
Shift-Row 54 32
Shift-Row 55 41389
Shift-Row 56 0

# Row 36 is the literal "-" placeholder row; keep it shifted too (value
# unchanged since every column already holds "-").
$ws.Range("D36").Value2 = "-"
$ws.Range("E36").Value2 = "-"
$ws.Range("F36").Value2 = "-"
$ws.Range("G36").Value2 = "-"
$ws.Range("H36").Value2 = "-"

# --- Column widths: the width-31 column (publish/format column) shifts
# left from F to E, and the trailing column reverts to the default 29.
$ws.Columns("D").ColumnWidth = 29
$ws.Columns("E").ColumnWidth = 31
$ws.Columns("F").ColumnWidth = 29
$ws.Columns("G").ColumnWidth = 29
$ws.Columns("H").ColumnWidth = 29

# --- Window size (cosmetic, best effort) ---------------------------------
$win = $excel.ActiveWindow
$win.Width = 20700
$win.Height = 9450
